$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.614.69'
$ws.Range('E2').Value = '  +3.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.699.77'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.98'
$ws.Range('E5').Value = '  +2.17%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3957'
$ws.Range('E7').Value = '  +2.10%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4055'
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '57.56'
$ws.Range('E9').Value = '  +17.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.561'
$ws.Range('E10').Value = '  +9.13%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08828'
$ws.Range('E12').Value = '  +2.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.352'
$ws.Range('E13').Value = '  +13.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.48'
$ws.Range('E14').Value = '  +4.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001338'
$ws.Range('E15').Value = '  +3.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.672'
$ws.Range('E16').Value = '  +7.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.693.94'
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '101.44'
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07054'
$ws.Range('E19').Value = '  +4.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.79'
$ws.Range('E20').Value = '  +4.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.945'
$ws.Range('E21').Value = '  +5.06%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.14'
$ws.Range('E23').Value = '  +2.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.593.36'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.970'
$ws.Range('E25').Value = '  +9.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.338'
$ws.Range('E26').Value = '  +1.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.55'
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.68'
$ws.Range('E28').Value = '  +2.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.233'
$ws.Range('E29').Value = '  +2.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.34'
$ws.Range('E30').Value = '  +3.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.624'
$ws.Range('E31').Value = '  +33.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.122'
$ws.Range('E32').Value = '  -3.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.880.19'
$ws.Range('E33').Value = '  +2.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.583'
$ws.Range('E34').Value = '  +18.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08594'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.26'
$ws.Range('E36').Value = '  +9.19%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.989'
$ws.Range('E37').Value = '  +1.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2767'
$ws.Range('E38').Value = '  +5.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.85'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02800'
$ws.Range('E40').Value = '  +11.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09085'
$ws.Range('E41').Value = '  +3.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.474'
$ws.Range('E42').Value = '  +2.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7784'
$ws.Range('E43').Value = '  +3.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7299'
$ws.Range('E44').Value = '  +4.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.80'
$ws.Range('E45').Value = '  +7.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.520'
$ws.Range('E46').Value = '  +6.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.200'
$ws.Range('E47').Value = '  +3.54%  '
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.70'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.292'
$ws.Range('E50').Value = '  +13.58%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00000000374'
$ws.Range('E51').Value = '  -0.28%  '
